$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 93
$lastRow  = 195
$newLastRow = 196

# --- 1. Snapshot the "before" values for the columns that move (D, J, K, M, P) ---
# for every row in [firstRow, lastRow] before any writes happen.
$oldD = @{}
$oldJ = @{}
$oldK = @{}
$oldM = @{}
$oldP = @{}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $oldD[$r] = $ws.Cells.Item($r, 4).Value2   # D: Fecha
    $oldJ[$r] = $ws.Cells.Item($r, 10).Value2  # J: Volumen
    $oldK[$r] = $ws.Cells.Item($r, 11).Value2  # K: Precio minimo
    $oldM[$r] = $ws.Cells.Item($r, 13).Value2  # M: Precio promedio ponderado
    $oldP[$r] = $ws.Cells.Item($r, 16).Value2  # P: Precio $/Kg
}

# Also snapshot the constant columns of row 195, which become the template
# for the brand-new row 196 (A,B,C,E,F,G,H,I,L,N,O,Q,R never change row to row).
$tplA = $ws.Cells.Item($lastRow, 1).Value2
$tplB = $ws.Cells.Item($lastRow, 2).Value2
$tplC = $ws.Cells.Item($lastRow, 3).Value2
$tplE = $ws.Cells.Item($lastRow, 5).Value2
$tplF = $ws.Cells.Item($lastRow, 6).Value2
$tplG = $ws.Cells.Item($lastRow, 7).Value2
$tplH = $ws.Cells.Item($lastRow, 8).Value2
$tplI = $ws.Cells.Item($lastRow, 9).Value2
$tplL = $ws.Cells.Item($lastRow, 12).Value2
$tplN = $ws.Cells.Item($lastRow, 14).Value2
$tplO = $ws.Cells.Item($lastRow, 15).Value2
$tplQ = $ws.Cells.Item($lastRow, 17).Value2
$tplR = $ws.Cells.Item($lastRow, 18).Value2
$dateFmt = $ws.Cells.Item($lastRow, 4).NumberFormat

# --- 2. Build the new row 196 (a copy of everything row 195 used to hold) ---
$ws.Cells.Item($newLastRow, 1).Value2 = $tplA
$ws.Cells.Item($newLastRow, 2).Value2 = $tplB
$ws.Cells.Item($newLastRow, 3).Value2 = $tplC
$ws.Cells.Item($newLastRow, 4).NumberFormat = $dateFmt
$ws.Cells.Item($newLastRow, 4).Value2 = $oldD[$lastRow]
$ws.Cells.Item($newLastRow, 5).Value2 = $tplE
$ws.Cells.Item($newLastRow, 6).Value2 = $tplF
$ws.Cells.Item($newLastRow, 7).Value2 = $tplG
$ws.Cells.Item($newLastRow, 8).Value2 = $tplH
$ws.Cells.Item($newLastRow, 9).Value2 = $tplI
$ws.Cells.Item($newLastRow, 10).Value2 = $oldJ[$lastRow]
$ws.Cells.Item($newLastRow, 11).Value2 = $oldK[$lastRow]
$ws.Cells.Item($newLastRow, 12).Value2 = $tplL
$ws.Cells.Item($newLastRow, 13).Value2 = $oldM[$lastRow]
$ws.Cells.Item($newLastRow, 14).Value2 = $tplN
$ws.Cells.Item($newLastRow, 15).Value2 = $tplO
$ws.Cells.Item($newLastRow, 16).Value2 = $oldP[$lastRow]
$ws.Cells.Item($newLastRow, 17).Value2 = $tplQ
$ws.Cells.Item($newLastRow, 18).Value2 = $tplR

# --- 3. Shift D/J/K/M/P down by one row: row r (94..195) takes row (r-1)'s old values ---
for ($r = $lastRow; $r -ge ($firstRow + 1); $r--) {
    $ws.Cells.Item($r, 4).Value2  = $oldD[$r - 1]
    $ws.Cells.Item($r, 10).Value2 = $oldJ[$r - 1]
    $ws.Cells.Item($r, 11).Value2 = $oldK[$r - 1]
    $ws.Cells.Item($r, 13).Value2 = $oldM[$r - 1]
    $ws.Cells.Item($r, 16).Value2 = $oldP[$r - 1]
}

# --- 4. Row 93 gets brand-new values ---
$ws.Cells.Item($firstRow, 4).Value2  = 44546
$ws.Cells.Item($firstRow, 10).Value2 = 2900
$ws.Cells.Item($firstRow, 11).Value2 = 400
$ws.Cells.Item($firstRow, 13).Value2 = 450
$ws.Cells.Item($firstRow, 16).Value2 = 900
